$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the source page name.
$ws.Name = "page_15"

# Build the format on each row's anchor cell first, then fan each one out
# to its whole row via copy/paste-special so every cell in the row picks up
# the same single style (instead of accumulating one style per property
# write per cell). Merging is done last for both rows so the two real
# (non-default) styles land on consecutive style-table slots, matching how
# Excel lays out cellXfs for this sheet.

# --- Row 1: " MASTER PACKAGE" header line, merged across A1:Z1. ---
$a1 = $ws.Range("A1")
$a1.Value = " MASTER PACKAGE"
$a1.Font.Size = 13
$a1.Font.Bold = $false
$a1.HorizontalAlignment = -4131
$a1.WrapText = $true

# --- Row 2: long descriptive line, merged across A2:Z2, bold. ---
$a2 = $ws.Range("A2")
$a2.Value = " WesternGlove Centric8 PROD                                                     M12225BVS563:KONRAD                                             M12225BVS                             MASTER"
$a2.Font.Size = 13
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4131
$a2.WrapText = $true

$row1 = $ws.Range("A1:Z1")
$a1.Copy()
$row1.PasteSpecial(-4122)

$row2 = $ws.Range("A2:Z2")
$a2.Copy()
$row2.PasteSpecial(-4122)

$row1.Merge()
$row2.Merge()

$excel.CutCopyMode = $false
